# Summon probability fix - 테스트용 확률 수정
$wb = $excel.ActiveWorkbook

# --- Rename third sheet ("Sheet2" -> "Follower") ---
$wsEquipment = $wb.Worksheets.Item("Equipment")
$wsSkills    = $wb.Worksheets.Item("Skills")
$wsFollower  = $wb.Worksheets.Item("Sheet2")
$wsFollower.Name = "Follower"

# --- Equipment sheet: selection moves from the whole-sheet range to H22 ---
$wsEquipment.Activate()
$wsEquipment.Range("H22").Select() | Out-Null

# --- Skills sheet: probability tweaks + selection move to F10 ---
$wsSkills.Range("C3").Value  = 3000
$wsSkills.Range("C4").Value  = 1996
$wsSkills.Range("C9").Value  = 3000
$wsSkills.Range("C10").Value = 1996

$wsSkills.Activate()
$wsSkills.Range("F10").Select() | Out-Null

# --- Follower sheet (formerly Sheet2): probability tweaks + selection/view move ---
$wsFollower.Range("C3").Value   = 1000
$wsFollower.Range("C4").Value   = 1000
$wsFollower.Range("C5").Value   = 1000
$wsFollower.Range("C6").Value   = 1000
$wsFollower.Range("C7").Value   = 999
$wsFollower.Range("C9").Value   = 1000
$wsFollower.Range("C10").Value  = 1000
$wsFollower.Range("C11").Value  = 1000
$wsFollower.Range("C12").Value  = 1000
$wsFollower.Range("C13").Value  = 999

$wsFollower.Activate()
$wsFollower.Range("F7").Select() | Out-Null
# View had scrolled so row 2 is the top visible row (topLeftCell = A2)
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
